# Added new test cases in Type Ahead service
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clean up stray / now-unused cells in the existing rows (2-6).
#    These used to hold a "STORE" value of PASS (column L) plus a
#    handful of structurally-empty cells that are no longer part of
#    the used range.
# ---------------------------------------------------------------------
$ws.Range("L2").ClearContents()

$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()

$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()

$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()

$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()

# ---------------------------------------------------------------------
# 2. New test cases (rows 7-12) for the Type Ahead service.
# ---------------------------------------------------------------------

# Row 7 - S1_TC_T6
$ws.Range("A7").Value2 = "S1_TC_T6"
$ws.Range("B7").Value2 = "Get Type Ahead by passing query,source and info values"
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value2 = "1PTYPEAHEAD"
$ws.Range("D7").Value2 = "/suggest"
$ws.Range("E7").Value2 = "GET"
$ws.Range("G7").Value2 = "?query=bio&source=wos&info=sports"
$ws.Range("H7").Value2 = "tmp"
$ws.Range("H7").ClearContents()
$ws.Range("J7").Value2 = "status=200||source=wos||suggestions.keyword=bio"
$ws.Rows.Item(7).RowHeight = 30

# Row 8 - S1_TC_T7
$ws.Range("A8").Value2 = "S1_TC_T7"
$ws.Range("B8").Value2 = "Get Type Ahead Suggestions for given query prefix, source and info"
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value2 = "1PTYPEAHEAD"
$ws.Range("D8").Value2 = "/suggest/ext/act"
$ws.Range("E8").Value2 = "GET"
$ws.Range("G8").Value2 = "?query=biology&source=wos&info=sports&size=1"
$ws.Range("H8").Value2 = "tmp"
$ws.Range("H8").ClearContents()
$ws.Range("J8").Value2 = "status=200||source=wos||suggestions.keyword=biology"
$ws.Rows.Item(8).RowHeight = 30

# Row 9 - S1_TC_T8 (healthcheck)
$ws.Range("A9").Value2 = "S1_TC_T8"
$ws.Range("B9").Value2 = "To verify HealthCheck"
$ws.Range("B9").WrapText = $true
$ws.Range("C9").Value2 = "1PTYPEAHEAD"
$ws.Range("D9").Value2 = "/healthcheck"
$ws.Range("E9").Value2 = "GET"
$ws.Range("H9").Value2 = "tmp"
$ws.Range("H9").ClearContents()
$ws.Range("J9").Value2 = "status=200"

# Row 10 - S1_TC_T9
$ws.Range("A10").Value2 = "S1_TC_T9"
$ws.Range("B10").Value2 = "Get Type Ahead Suggestions for given query prefix, multiple source values and info"
$ws.Range("B10").WrapText = $true
$ws.Range("C10").Value2 = "1PTYPEAHEAD"
$ws.Range("D10").Value2 = "/suggest/ext/act"
$ws.Range("E10").Value2 = "GET"
$ws.Range("G10").Value2 = "?query=biology&source=wos&source=categories&info=sports&size=1"
$ws.Range("H10").Value2 = "tmp"
$ws.Range("H10").ClearContents()
$ws.Range("J10").Value2 = "status=200||source=wos||source=categories||suggestions.keyword=biology"
$ws.Rows.Item(10).RowHeight = 30

# Row 11 - S1_TC_T10
$ws.Range("A11").Value2 = "S1_TC_T10"
$ws.Range("B11").Value2 = "Get Type Ahead by passing multiple sources and info values."
$ws.Range("B11").WrapText = $true
$ws.Range("C11").Value2 = "1PTYPEAHEAD"
$ws.Range("D11").Value2 = "/suggest"
$ws.Range("E11").Value2 = "GET"
$ws.Range("G11").Value2 = "?query=biology&source=wos&source=categories&info=sports&size=1"
$ws.Range("H11").Value2 = "tmp"
$ws.Range("H11").ClearContents()
$ws.Range("J11").Value2 = "status=200||source=wos||source=categories||suggestions.keyword=biology"
$ws.Rows.Item(11).RowHeight = 30

# Row 12 - S1_TC_T11 (invalid query params - size)
$ws.Range("A12").Value2 = "S1_TC_T11"
$ws.Range("B12").Value2 = "To validate Type Ahead response for invalid query params - size"
$ws.Range("B12").WrapText = $true
$ws.Range("C12").Value2 = "1PTYPEAHEAD"
$ws.Range("D12").Value2 = "/suggest/ext/act"
$ws.Range("E12").Value2 = "GET"
$ws.Range("G12").Value2 = "?query=biology&source=wos&source=categories&info=sports&size=a"
$ws.Range("H12").Value2 = "tmp"
$ws.Range("H12").ClearContents()
$ws.Range("J12").Value2 = "status=404"
$ws.Rows.Item(12).RowHeight = 30

# ---------------------------------------------------------------------
# 4. Update the view: selection now spans L2:L12, and the sheet is
#    scrolled so column E is leftmost.
# ---------------------------------------------------------------------
$ws.Range("L2:L12").Select()
$excel.ActiveWindow.ScrollColumn = 5
